$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "1.000", "298.75") that
# Excel would otherwise auto-convert to real numbers on assignment. Force the
# cell to Text first, assign, then restore the default "Normal" style so the
# saved cell matches the original (unstyled) inline-string cells.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "23.464.75"
$ws.Range("E2").Value = "  -1.72%  "
Set-TextValue $ws.Range("D3") "1.648.12"
$ws.Range("E3").Value = "  -0.51%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.46%  "
Set-TextValue $ws.Range("D5") "1.001"
$ws.Range("E5").Value = "  +0.24%  "
Set-TextValue $ws.Range("D6") "298.75"
$ws.Range("E6").Value = "  -1.84%  "
Set-TextValue $ws.Range("D7") "0.3782"
$ws.Range("E7").Value = "  -1.24%  "
Set-TextValue $ws.Range("D8") "0.3551"
$ws.Range("E8").Value = "  -1.98%  "
Set-TextValue $ws.Range("D9") "49.90"
$ws.Range("E9").Value = "  -2.75%  "
Set-TextValue $ws.Range("D10") "0.08094"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -2.87%  "
Set-TextValue $ws.Range("D12") "1.000"
$ws.Range("E12").Value = "  +0.48%  "
Set-TextValue $ws.Range("D13") "22.09"
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("E14").Value = "  -2.57%  "
Set-TextValue $ws.Range("D15") "7.358"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("E16").Value = "  -3.36%  "
Set-TextValue $ws.Range("D17") "1.649.35"
$ws.Range("E17").Value = "  -0.29%  "
Set-TextValue $ws.Range("D18") "97.26"
$ws.Range("E18").Value = "  -0.52%  "
Set-TextValue $ws.Range("D19") "0.06939"
$ws.Range("E19").Value = "  -0.62%  "
Set-TextValue $ws.Range("D20") "6.751"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("E22").Value = "  +0.11%  "
Set-TextValue $ws.Range("D23") "12.44"
$ws.Range("E23").Value = "  -1.79%  "
Set-TextValue $ws.Range("D24") "23.466.95"
$ws.Range("E24").Value = "  -1.65%  "
Set-TextValue $ws.Range("D25") "2.493"
$ws.Range("E25").Value = "  -1.77%  "
Set-TextValue $ws.Range("D26") "2.905"
$ws.Range("E26").Value = "  -5.85%  "
Set-TextValue $ws.Range("D27") "20.91"
$ws.Range("E27").Value = "  -2.04%  "
Set-TextValue $ws.Range("D28") "152.52"
$ws.Range("E28").Value = "  +0.86%  "
Set-TextValue $ws.Range("D29") "5.201"
$ws.Range("E29").Value = "  -0.59%  "
Set-TextValue $ws.Range("D30") "132.71"
$ws.Range("E30").Value = "  -1.56%  "
Set-TextValue $ws.Range("D31") "1.835.20"
Set-TextValue $ws.Range("D32") "6.928"
$ws.Range("E32").Value = "  +0.34%  "
Set-TextValue $ws.Range("D33") "2.127"
$ws.Range("E33").Value = "  +0.99%  "
Set-TextValue $ws.Range("D34") "11.44"
$ws.Range("E34").Value = "  -5.43%  "
Set-TextValue $ws.Range("D35") "0.9958"
$ws.Range("E35").Value = "  -8.39%  "
Set-TextValue $ws.Range("D36") "0.02719"
$ws.Range("E36").Value = "  -4.49%  "
Set-TextValue $ws.Range("D37") "0.08709"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D38") "0.2431"
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D39") "5.937"
$ws.Range("E39").Value = "  -4.03%  "
Set-TextValue $ws.Range("D40") "13.07"
$ws.Range("E40").Value = "  +1.14%  "
Set-TextValue $ws.Range("D41") "0.06780"
$ws.Range("E41").Value = "  -4.15%  "
Set-TextValue $ws.Range("D42") "0.6891"
$ws.Range("E42").Value = "  -2.99%  "
Set-TextValue $ws.Range("D43") "1.303"
$ws.Range("E43").Value = "  -3.01%  "
Set-TextValue $ws.Range("D44") "15.56"
$ws.Range("E44").Value = "  -1.97%  "
Set-TextValue $ws.Range("D45") "0.9998"
$ws.Range("E45").Value = "  +0.20%  "
Set-TextValue $ws.Range("D46") "0.6361"
$ws.Range("E46").Value = "  -3.38%  "
Set-TextValue $ws.Range("D47") "2.258"
$ws.Range("E47").Value = "  -3.67%  "
Set-TextValue $ws.Range("D48") "3.904"
$ws.Range("E48").Value = "  -1.88%  "
Set-TextValue $ws.Range("D49") "0.07726"
$ws.Range("E49").Value = "  -3.48%  "
Set-TextValue $ws.Range("D50") "126.74"
$ws.Range("E50").Value = "  -2.03%  "
Set-TextValue $ws.Range("D51") "1.151"
$ws.Range("E51").Value = "  -4.26%  "
